$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename Importe -> Importe_ARS, add Importe_USD column ---
$ws.Range("B1").Value = "Importe_ARS"
$ws.Range("C1").Value = "Importe_USD"

# Give the new header cell (C1) the same formatting as the existing header cells
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- Data rows: amounts become formatted text strings, new USD column added ---
$ws.Range("B2").Value = "45.666,00"
$ws.Range("C2").Value = "0,00"

$ws.Range("B3").Value = "98.777,00"
$ws.Range("C3").Value = "0,00"

$ws.Range("B4").Value = "144.443,00"
$ws.Range("C4").Value = "0,00"
